$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row -> new F value }
$updates = @{
    "展览" = @{
        2  = 76
        5  = 31
        6  = 21
        7  = 254
        8  = 14174
        9  = 136
        10 = 104
        11 = 5690
        13 = 61
        15 = 55
        16 = 1229
        18 = 75
        19 = 171
        20 = 773
        21 = 2917
        22 = 50
        23 = 10477
        24 = 1192
        25 = 47
        26 = 69
        27 = 3721
        28 = 240
    }
    "全部类型" = @{
        2  = 76
        6  = 31
        7  = 21
        8  = 254
        9  = 14174
        10 = 136
        11 = 104
        12 = 5690
        14 = 61
        16 = 55
        17 = 1229
        19 = 75
        20 = 171
        21 = 773
        22 = 2917
        23 = 50
        25 = 10477
        26 = 1192
        27 = 47
        28 = 69
        29 = 3721
        30 = 240
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($r in $rows.Keys) {
        $ws.Range("F$r").Value = $rows[$r]
    }
}
